# Mise à jour de l'application
# Adds 7 new training-log rows (273-279, all dated 2025-09-04 / serial 45904)
# to the bottom of the "Feuil1" worksheet, extends the shared charge formula
# and refreshes the current selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# --- 1. Stamp formatting for the new rows by cloning existing template rows ---
# Row 272 has an empty "Localisation douleur" (G) cell -> style pattern to
# reuse for new rows whose G column stays blank.
# Row 271 has a text value in G -> style pattern to reuse for new rows that
# also carry a text value in G.
$ws.Range("A272:I272").Copy() | Out-Null
$ws.Range("A273:I273").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A275:I277").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("A271:I271").Copy() | Out-Null
$ws.Range("A274:I274").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A278:I279").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# --- 2. New row data --------------------------------------------------------
# Date | Name | Volume | Intensite | Fatigue | Douleur | Localisation | Plaisir
$rows = @(
    @{ Row = 273; Date = 45904; Name = "Amir Etien";       Volume = 50; Intensite = 4; Fatigue = 7; Douleur = 0; Localisation = $null;           Plaisir = 3 },
    @{ Row = 274; Date = 45904; Name = "Naim Dhib";        Volume = 50; Intensite = 3; Fatigue = 5; Douleur = 3; Localisation = "Courbatures ";  Plaisir = 5 },
    @{ Row = 275; Date = 45904; Name = "Ilyes Boughanmi";  Volume = 50; Intensite = 5; Fatigue = 5; Douleur = 0; Localisation = $null;           Plaisir = 10 },
    @{ Row = 276; Date = 45904; Name = "Omar Benyounes";   Volume = 50; Intensite = 5; Fatigue = 6; Douleur = 0; Localisation = $null;           Plaisir = 4 },
    @{ Row = 277; Date = 45904; Name = "Yanis Berrached";  Volume = 50; Intensite = 3; Fatigue = 8; Douleur = 0; Localisation = $null;           Plaisir = 6 },
    @{ Row = 278; Date = 45904; Name = "Emmanuel Valey";   Volume = 50; Intensite = 1; Fatigue = 2; Douleur = 2; Localisation = "Adducteur ";    Plaisir = 7 },
    @{ Row = 279; Date = 45904; Name = "Karahali Souaré";  Volume = 50; Intensite = 4; Fatigue = 6; Douleur = 7; Localisation = "Cheville ";     Plaisir = 7 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $r.Name
    $ws.Cells.Item($row, 3).Value = $r.Volume
    $ws.Cells.Item($row, 4).Value = $r.Intensite
    $ws.Cells.Item($row, 5).Value = $r.Fatigue
    $ws.Cells.Item($row, 6).Value = $r.Douleur
    if ($r.Localisation) {
        $ws.Cells.Item($row, 7).Value = $r.Localisation
    }
    $ws.Cells.Item($row, 8).Value = $r.Plaisir
    $ws.Cells.Item($row, 9).Formula = "=C$row*D$row"
}

# --- 3. Refresh the view: scroll position & active selection ---------------
$ws.Range("K274").Select()

